$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data range (A1:E2) before writing the new layout
$ws.Range("A1:E2").Clear()

# New data: bus number paired with its reference voltage, one bus per row
$data = @(
    @(1, 1.05),
    @(2, 1.0449999999999999),
    @(3, 1.01),
    @(10, 1.05),
    @(12, 1.05)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Update selection to match the new data extent
$ws.Range("A1:B5").Select()
